$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Delete the duplicate "devices (2)" worksheet entirely.
$oldDevicesCopy = $wb.Worksheets.Item("devices (2)")
$oldDevicesCopy.Delete()

# 2. Update the data on the "devices" sheet: the sign-in device row now
#    references a person (Brian) and browser (SAFARI) instead of
#    Avner / mobileOS.
$devices = $wb.Worksheets.Item("devices")
$devices.Range("I2").Value = "Brian"
$devices.Range("J2").Value = "SAFARI"

# 3. Move "devices" to be the first sheet in the workbook and make it the
#    active/selected tab.
$devices.Move($wb.Worksheets.Item(1))
$devices = $wb.Worksheets.Item("devices")
$devices.Activate()
$devices.Range("I2").Select()

# 4. Update the selection remembered on the "signIn" sheet.
$signIn = $wb.Worksheets.Item("signIn")
$signIn.Activate()
$signIn.Range("B17").Select()

# 5. Restore focus back onto the "devices" tab, which is the active sheet
#    in the final workbook.
$devices.Activate()
